# "Generate Report for Handback"
#
# The handback report workbook records, per localized file, the timestamps
# of the latest handoff/handback xliff round-trip. This script simulates a
# fresh report generation run that refreshes those timestamps for the
# "f0820fe7-175a-4264-a194-9d9298b01209.md" file (row 3 in every sheet),
# while leaving the "57148685-a294-4851-acbf-0eb1ac145cfc.md" file (row 2)
# untouched, since it had already been handed back earlier and was not
# reprocessed in this run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
# "Latest HO Xliff Generate Date" column (G) for the f0820fe7 file row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-30 16:57:23"

# --- zh-cn sheet -----------------------------------------------------------
# New handoff/handback round trip datetimes for the f0820fe7 file row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-30 16:57:18"
$wsZhCn.Range("K3").Value = "2016-08-30 16:57:36"

# --- de-de sheet -------------------------------------------------------
# New handoff/handback round trip datetimes for the f0820fe7 file row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-30 16:57:23"
$wsDeDe.Range("K3").Value = "2016-08-30 16:57:44"
